$d = $word.ActiveDocument

# 1) "What are the number ... College of University? <EN-DASH> Drop Down"
#    -> "Total number of institutes affiliated by A.I.C.T.E.?"
$enDash = [char]0x2013
$old1 = "What are the number of grievances filed, resolved in College of University? " + $enDash + " Drop Down"
$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false,
    "Total number of institutes affiliated by A.I.C.T.E.?", 2) | Out-Null

# 2) "How many grievances are pending for the university/college?"
#    -> "Total number of grievances reported?"
$d.Content.Find.Execute("How many grievances are pending for the university/college?", $true, $false, $false, $false, $false, $true, 1, $false,
    "Total number of grievances reported?", 2) | Out-Null

# 3) TAB + "Total number of institutes affiliated by A.I.C.T.E.?"
#    -> "Total number of addressed grievances?" (tab removed)
$tab = [char]9
$d.Content.Find.Execute($tab + "Total number of institutes affiliated by A.I.C.T.E.?", $true, $false, $false, $false, $false, $true, 1, $false,
    "Total number of addressed grievances?", 2) | Out-Null

# 4) TAB + "Total number of grievances reported?"
#    -> "Total number of In Action Grievances" (tab removed)
$d.Content.Find.Execute($tab + "Total number of grievances reported?", $true, $false, $false, $false, $false, $true, 1, $false,
    "Total number of In Action Grievances", 2) | Out-Null

# 5) TAB + "Total number of addressed grievances?"
#    -> "Total number of Delayed Grievances" (tab removed)
$d.Content.Find.Execute($tab + "Total number of addressed grievances?", $true, $false, $false, $false, $false, $true, 1, $false,
    "Total number of Delayed Grievances", 2) | Out-Null

# 6/7/8) Collapse the final three paragraphs:
#   - "How many grievances were lodged in the last 6 months?" (has lastRenderedPageBreak)
#   - "Which college/university had the maximum grievances lodged, resolved in the last 6 months" (removed entirely)
#   - the trailing ind-left=360 paragraph that held the _GoBack bookmark
# into:
#   - a paragraph reading "Total number of Re-Open Grievances" (keeps lastRenderedPageBreak) with the
#     _GoBack bookmark now anchored at its end
#   - the trailing ind-left=360 paragraph (now without the bookmark)
$pCount = $d.Paragraphs.Count
$pStart = $d.Paragraphs.Item($pCount - 2)
$pEnd = $d.Paragraphs.Item($pCount)
$tailRange = $d.Range($pStart.Range.Start, $pEnd.Range.End)

$tailXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
  '<w:body>' +
  '<w:p w14:paraId="76DFA852" w14:textId="7959062C" w:rsidR="0082531E" w:rsidRPr="00127AD3" w:rsidRDefault="0082531E" w:rsidP="0082531E">' +
  '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr>' +
  '<w:r><w:lastRenderedPageBreak/><w:t>Total number of Re-Open Grievances</w:t></w:r>' +
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>' +
  '<w:p w14:paraId="7DDF79B1" w14:textId="77777777" w:rsidR="00ED026D" w:rsidRPr="00F5691A" w:rsidRDefault="00ED026D" w:rsidP="00107ACE">' +
  '<w:pPr><w:ind w:left="360"/></w:pPr></w:p>' +
  '</w:body></w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'

$tailRange.InsertXML($tailXml)
